# --------------------------------------------------------------------------
# "Update countries & provincias Spain" - refresh of the COVID-19 "Pais" sheet
# --------------------------------------------------------------------------
# 1) The "last updated" banner in A1 moves from 14:00 to 15:17.
# 2) Three country-ranking pairs swap places (their totals overtook one
#    another since the previous refresh), so the country names in column A
#    need to be swapped for those row pairs:
#       row 39/40   : Israel  / Kuwait   -> Kuwait / Israel
#       row 81/82   : Bulgaria / Bosnia y Herzegovina -> Bosnia y Herzegovina / Bulgaria
#       row 84/85   : Senegal / Madagascar -> Madagascar / Senegal
# 3) Updated per-country statistics (B:H = Casos totales, Nuevos casos,
#    Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes).
# --------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Refreshed timestamp
$ws.Range("A1").Value = "Datos actualizados a 28 de Julio de 2020 a las 15:17"

# 2) Country-name swaps for the rows whose rank order flipped
$ws.Range("A39").Value = "Kuwait"
$ws.Range("A40").Value = "Israel"
$ws.Range("A81").Value = "Bosnia y Herzegovina"
$ws.Range("A82").Value = "Bulgaria"
$ws.Range("A84").Value = "Madagascar"
$ws.Range("A85").Value = "Senegal"

# 3) Updated statistics, keyed by row number -> column letter -> new value.
#    Only the listed columns change for each row; every other cell (e.g. the
#    "Casos criticos" column F, which stays 0 everywhere) is left untouched.
$rowUpdates = [ordered]@{
    4 = @{ "B" = 4434185; "C" = 775; "D" = 2137959; "E" = 2145726; "G" = 56; "H" = 150500 }
    6 = @{ "B" = 1493904; "C" = 11401; "D" = 957044; "E" = 503323; "G" = 89; "H" = 33537 }
    16 = @{ "B" = 270831; "C" = 1897; "D" = 225624; "E" = 42418; "G" = 29; "H" = 2789 }
    36 = @{ "B" = 67366; "C" = 115; "D" = 60669; "E" = 6154; "G" = 5; "H" = 543 }
    39 = @{ "B" = 65149; "C" = 770; "D" = 55681; "E" = 9026; "G" = 4; "H" = 442 }
    40 = @{ "B" = 64649; "C" = 664; "D" = 32117; "E" = 32052; "G" = 6; "H" = 480 }
    44 = @{ "B" = 53374; "C" = 223; "G" = 4; "H" = 6145 }
    46 = @{ "B" = 50410; "C" = 111; "D" = 35626; "E" = 13062; "G" = 3; "H" = 1722 }
    55 = @{ "B" = 34609; "C" = 132; "E" = 1731 }
    58 = @{ "B" = 30858; "C" = 412; "D" = 23873; "E" = 6555; "G" = 7; "H" = 430 }
    67 = @{ "B" = 19063; "C" = 311; "D" = 13875; "E" = 5139; "G" = 1; "H" = 49 }
    78 = @{ "B" = 13577; "C" = 30; "D" = 12451; "E" = 513 }
    80 = @{ "B" = 10938; "C" = 317; "E" = 7108 }
    81 = @{ "B" = 10766; "C" = 268; "D" = 5220; "E" = 5249; "G" = 3; "H" = 297 }
    82 = @{ "B" = 10621; "D" = 5585; "E" = 4689; "H" = 347 }
    83 = @{ "B" = 10315; "C" = 102; "D" = 5663; "E" = 4181; "G" = 5; "H" = 471 }
    84 = @{ "B" = 10104; "C" = 414; "D" = 6613; "E" = 3398; "G" = 2; "H" = 93 }
    85 = @{ "B" = 9805; "C" = 41; "D" = 6591; "E" = 3016; "G" = 4; "H" = 198 }
    88 = @{ "B" = 8873; "C" = 29; "D" = 5930; "E" = 2735 }
    98 = @{ "B" = 4923; "C" = 42; "D" = 4034; "E" = 749; "G" = 1; "H" = 140 }
    116 = @{ "B" = 2807; "C" = 2; "E" = 500 }
    120 = @{ "B" = 2520; "C" = 7; "D" = 1919; "E" = 477 }
    130 = @{ "B" = 1857; "C" = 3; "E" = 24 }
    156 = @{ "B" = 708; "C" = 7; "E" = 34 }
    163 = @{ "B" = 438; "C" = 7; "D" = 369; "E" = 69 }
    176 = @{ "B" = 220; "C" = 6; "E" = 32 }
    179 = @{ "B" = 186; "C" = 1; "E" = 6 }
}

foreach ($row in $rowUpdates.Keys) {
    $cols = $rowUpdates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
